$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.084.53"
$ws.Range("E2").Value = "  +4.60%  "
$ws.Range("D3").Value = "2.262.33"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("D5").Value = "'253.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'0.640"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("D7").Value = "'71.76"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.94%  "
$ws.Range("D8").Value = "'0.674"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +17.98%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'40.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.23%  "
$ws.Range("D11").Value = "'0.0975"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.99%  "
$ws.Range("D12").Value = "'59.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "'7.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.49%  "
$ws.Range("D14").Value = "'0.104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "2.607.65"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.887"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'14.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "2.260.43"
$ws.Range("E18").Value = "  +3.83%  "
$ws.Range("D19").Value = "42.996.02"
$ws.Range("E19").Value = "  +4.26%  "
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").Value = "'73.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("D23").Value = "'236.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.74%  "
$ws.Range("D25").Value = "'3.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").Value = "'11.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "'2.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D31").Value = "'167.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "'21.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").Value = "'0.129"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.84%  "
$ws.Range("D34").Value = "'6.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.94%  "
$ws.Range("D35").Value = "'0.0778"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.25%  "
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").Value = "'28.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.12%  "
$ws.Range("D38").Value = "'4.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.63%  "
$ws.Range("D39").Value = "'4.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").Value = "'0.0322"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.88%  "
$ws.Range("D41").Value = "'2.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.20%  "
$ws.Range("D42").Value = "'5.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("D43").Value = "'12.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("D44").Value = "'64.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").Value = "'4.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("D47").Value = "'8.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("D49").Value = "'1.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").Value = "'1.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").Value = "  +2.32%  "
